$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.988.45"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "1.885.80"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.69"
$ws.Range("E5").Value = "  -2.31%  "

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4586"
$ws.Range("E7").Value = "  -2.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4070"
$ws.Range("E8").Value = "  +0.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.45"
$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07968"
$ws.Range("E10").Value = "  -2.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9915"
$ws.Range("E11").Value = "  -2.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.67"
$ws.Range("E12").Value = "  -3.10%  "

$ws.Range("D13").Value = "1.883.25"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.902"
$ws.Range("E14").Value = "  -3.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.059"
$ws.Range("E15").Value = "  -3.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.45"
$ws.Range("E17").Value = "  -3.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001025"
$ws.Range("E18").Value = "  -2.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06550"
$ws.Range("E19").Value = "  -0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.43"
$ws.Range("E20").Value = "  -1.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "29.003.12"
$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.413"
$ws.Range("E23").Value = "  -2.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.27"
$ws.Range("E24").Value = "  +1.29%  "

$ws.Range("E25").Value = "  -2.70%  "

$ws.Range("D26").Value = "2.102.51"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.76"
$ws.Range("E27").Value = "  -2.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.57"
$ws.Range("E28").Value = "  -2.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.100"
$ws.Range("E29").Value = "  -2.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.403"
$ws.Range("E30").Value = "  -2.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.72"
$ws.Range("E31").Value = "  -2.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9965"
$ws.Range("E32").Value = "  -1.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09326"
$ws.Range("E33").Value = "  -2.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.405"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.280"
$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06047"
$ws.Range("E37").Value = "  -2.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02215"
$ws.Range("E38").Value = "  -3.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.265"
$ws.Range("E39").Value = "  -4.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.181"
$ws.Range("E40").Value = "  -1.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9996"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5782"
$ws.Range("E42").Value = "  -3.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1818"
$ws.Range("E43").Value = "  -4.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.10"
$ws.Range("E44").Value = "  -3.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.253"
$ws.Range("E45").Value = "  -1.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07523"
$ws.Range("E46").Value = "  +3.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.269"
$ws.Range("E47").Value = "  +7.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.97"
$ws.Range("E48").Value = "  -2.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5452"
$ws.Range("E49").Value = "  -2.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.894"
$ws.Range("E50").Value = "  -4.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.15"
$ws.Range("E51").Value = "  -1.49%  "
